$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-01-30", "14:45:02", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:45:13", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:45:23", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:45:34", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "14:45:53", "14:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED"),
    @("2026-01-30", "14:45:54", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$defaultStyle = $ws.Cells.Item(1, 1).Style

$startRow = 116
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a plain-looking date string (e.g. "2026-01-30"); force it
    # to be stored as text (like the rest of the column) instead of letting
    # Excel auto-convert it into a date serial number.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.Style = $defaultStyle

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
